$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the "Meta description" paragraph (currently paragraph 2,
#    right after the title heading) and remember its bold-run
#    formatted text so we can reuse it (same bold-run shape) further
#    down in the document. Do this - and the insertion below - before
#    deleting the paragraph, so the saved range/formatting stays
#    valid (deleting first would shift every later offset).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaBoldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 16)  # "Meta description"
$boldFormattedText = $metaBoldRange.FormattedText

# ------------------------------------------------------------------
# 2. Find the paragraph right before the last one (the one that
#    currently reads "Lack of variety in base game symbols") and add
#    a new paragraph after it, reusing the bold-run formatting copied
#    above, then overwrite its text with the new title line.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$beforeLastPara = $d.Paragraphs.Item($lastIndex - 1)
$beforeLastPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($lastIndex)
$newPara.Style = "Normal"

$insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertPoint.FormattedText = $boldFormattedText

$newTitleRange = $newPara.Range
$newTitleRange.MoveEnd(1, -1) | Out-Null
$newTitleRange.Text = "Play Christmas Gold Digger Slot Free | Festive 5x3 Grid Game"

# ------------------------------------------------------------------
# 3. Now that the copy is complete, delete the original "Meta
#    description" paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 4. Replace the final paragraph's "Dear DALLE..." text with the new
#    meta-description copy, keeping its existing italic run
#    formatting untouched.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Dear DALLE, I need a feature image for the online slot game " + [char]34 + "Christmas Gold Digger" + [char]34 + ". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The main focus of the image should be on the Maya warrior, with elements of Christmas and gold mining included in the background. The image should look festive and fun, with bright colors and playful details that will appeal to players. Please include the game's title " + [char]34 + "Christmas Gold Digger" + [char]34 + " in the image, as well as any other elements that you think will help players identify and enjoy this game. Thank you!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Christmas Gold Digger online slot. Play for free with festive graphics, special gold hold feature, and autoplay with win/loss limits.",
    2
) | Out-Null
